# Insert a new record row at the top of the Coliflor (Macroferia Regional de
# Talca) weekly block.  The sheet already has rows 2..396 of data; row 376 is
# the first row of this particular market/product block.  We insert a new
# row at position 376, which pushes the existing rows 376..396 down to
# 377..397, and then populate the freshly inserted row with the new weekly
# observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 376 (shifts 376-396 -> 377-397)
$ws.Rows.Item(376).Insert()

# Populate the new row 376 with the new weekly data point
$ws.Range("A376").Value = 5
$ws.Range("B376").Value = "Macroferia Regional de Talca"
$ws.Range("C376").Value = "Maule"
$ws.Range("D376").Value = 45013
$ws.Range("E376").Value = 7
$ws.Range("F376").Value = 100112008
$ws.Range("G376").Value = "Coliflor"
$ws.Range("H376").Value = "Sin especificar"
$ws.Range("I376").Value = "Primera"
$ws.Range("J376").Value = 4000
$ws.Range("K376").Value = 800
$ws.Range("L376").Value = 900
$ws.Range("M376").Value = 850
$ws.Range("N376").Value = "`$/unidad"
$ws.Range("O376").Value = "Región del Maule"
$ws.Range("P376").Value = 850
$ws.Range("Q376").Value = 1
$ws.Range("R376").Value = "Hortaliza"
